# Applies the RS-RDR:15-15 resourcesResponse schema update.
# NOTE: this runtime's Range.Find.Execute operates over the whole document
# content regardless of which Range/Selection invoked it, so any edit whose
# search text is not unique in the document is done via direct Cell/Range
# Text assignment (which IS properly scoped) instead of Find/Replace.

$d = $word.ActiveDocument
$vt = [char]11   # manual line break within a run -> <w:br/>

# ---------------------------------------------------------------------
# 1) Title heading
# ---------------------------------------------------------------------
$d.Paragraphs.Item(1).Range.Text = "Objet RS-RDR:15-15:resourcesResponse"

# ---------------------------------------------------------------------
# Table 1 : "Objet resourcesResponse" table (caseId / RSDDRId / resourceRequest / mobilizedResources)
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(1)

# 2) caseId example: dash -> dot
$t1.Cell(2,6).Range.Text = "fr.health.samu440.DRFR15DDXAAJJJ0000"

# 3) mobilizedResources "Champ correspondant"
$t1.Cell(5,2).Range.Text = "Ressource engagée"

# ---------------------------------------------------------------------
# Table 2 : "Type request" table (dateTime / answer / deadline / freetext)
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

# 4) answer.Format: add ENUM break
$t2.Cell(3,3).Range.Text = "string" + $vt + "(ENUM : OUI, NON, PARTIEL, DIFFERE)"

# 5) answer.Description: complete the enumeration
$t2.Cell(3,5).Range.Text = "oui / non / oui partiel / différé"

# ---------------------------------------------------------------------
# Table 3 : "Type resource" table (19 rows)
# ---------------------------------------------------------------------
$t3 = $d.Tables.Item(3)

# 6) resourceID.Description + example
$t3.Cell(5,5).Range.Text = "ID unique de la ressource engagée partagée = aux champs {orgID}.R.{ownerID}"
$t3.Cell(5,6).Range.Text = "fr.health.samu440.R.123456"

# 7) orgID example
$t3.Cell(6,6).Range.Text = "fr.health.samu440"

# 8) "type" row -> "resourceType" row, with updated format/description/example
$t3.Cell(7,1).Range.Text = "resourceType"
$t3.Cell(7,2).Range.Text = "Type de ressource"
$t3.Cell(7,3).Range.Text = "string" + $vt + "(NOMENCLATURE: CISU-TYPE_MOYEN)"
$t3.Cell(7,5).Range.Text = "Type de ressource mobilisée : Smur, Hospitaliers (hors Smur), Professionnels Libéraux, Ambulanciers privés (Transporteurs Sanitaires Urgent), etc."
$t3.Cell(7,6).Range.Text = "SMUR"

# 9) New row "vehiculeType" inserted right after the resourceType row (row 8)
$t3.Rows.Add($t3.Rows.Item(8)) | Out-Null
$t3.Cell(8,1).Range.Text = "vehiculeType"
$t3.Cell(8,2).Range.Text = "Type de vecteur"
$t3.Cell(8,3).Range.Text = "string" + $vt + "(NOMENCLATURE: CISU-TYPE_VECTEUR)"
$t3.Cell(8,4).Range.Text = "0..1"
$t3.Cell(8,5).Range.Text = "Type de vecteur mobilisé : Véhicule Léger Médicalisé, Ambulance de réanimation, Ambulance de réanimation Bariatrique, Ambulance de réanimation Pédiatrique, etc."
$t3.Cell(8,6).Range.Text = "VLM"

# 10) Remove the trailing "freetext" row (now last row, after the insert above)
$t3.Rows.Item($t3.Rows.Count).Delete()

# ---------------------------------------------------------------------
# Table 4 : "Type team" table (type / name)
# ---------------------------------------------------------------------
$t4 = $d.Tables.Item(4)

# 11) type.Format ENUM accent fix
$t4.Cell(2,3).Range.Text = "string" + $vt + "(ENUM : Medicale, Paramedicale)"

# ---------------------------------------------------------------------
# Table 5 : "Type state" table (dateTime / status / availability)
# ---------------------------------------------------------------------
$t5 = $d.Tables.Item(5)

# 12) status.Format ENUM accent fix
$t5.Cell(3,3).Range.Text = "string" + $vt + "(ENUM : Alerte, Parti, Arrivee sur les lieux, Transport destination, Arrivee destination, Fin de medicalisation , Quitte destination, Retour base, Rentree Base)"

# 13) status example
$t5.Cell(3,6).Range.Text = "ENUM : Alerte, Parti, Arrivee Sur Les Lieux, Transport Destination, Arrivee Destination, Fin De Medicalisation , Quitte Destination, Retour Base, Rentree Base"

# 14) availability example
$t5.Cell(4,6).Range.Text = "ENUM : Disponible, Indisponible, Inconnu"

Write-Output "RS-RDR schema update applied."
